$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire row 4 (to mirror the selection change seen in the diff)
$ws.Rows(4).Select()

# Delete the selected row (this is row 4, "Universal Music Plaza Stage"),
# which shifts all subsequent rows up by one.
$ws.Rows(4).Delete()
